# Ion card - new registry (Turgenev)
#
# Adds a new registry row (№3: И.С. Тургенев - "Отцы и дети") to the
# loan/registry sheet, and highlights the header row in yellow to call
# it out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- New row 4: pick up the same bordered formatting already used by
#     the rest of the table before writing the new values into it ---
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial($xlPasteFormats)

$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial($xlPasteFormats)

$ws.Range("E1").Copy()
$ws.Range("C4:E4").PasteSpecial($xlPasteFormats)

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 44498
$ws.Range("B4").NumberFormat = "d-mmm"
$ws.Range("C4").Value = "И.С. Тургенев"
$ws.Range("D4").Value = "Отцы и дети"
$ws.Range("E4").Value = "стр.45"

# --- Header row: highlight with a yellow fill ---
$ws.Range("A1:E1").Interior.Color = 65535

# --- Selection now rests on the header row ---
$ws.Range("A1:E1").Select() | Out-Null

$excel.CutCopyMode = $false | Out-Null
